$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C32").Value = 8305
$ws.Range("C33").Value = 8274
$ws.Range("C34").Value = 8267
$ws.Range("C35:C76").Value = 7660
$ws.Range("C77:C88").Value = 7318
$ws.Range("C89:C139").Value = 7293
